$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "72.398.50"
Set-TextValue "E2" "  +4.32%  "

Set-TextValue "D3" "4.038.30"
Set-TextValue "E3" "  +3.21%  "

Set-TextValue "E4" "  +0.14%  "

Set-TextValue "D5" "518.34"
Set-TextValue "E5" "  -2.20%  "

Set-TextValue "D6" "147.28"
Set-TextValue "E6" "  +2.09%  "

Set-TextValue "D7" "0.735"
Set-TextValue "E7" "  +19.88%  "

Set-TextValue "E8" "  +0.13%  "

Set-TextValue "D9" "0.771"
Set-TextValue "E9" "  +7.34%  "

Set-TextValue "E10" "  +0.73%  "

Set-TextValue "D11" "0.0000329"
Set-TextValue "E11" "  -2.26%  "

Set-TextValue "D12" "47.26"
Set-TextValue "E12" "  +12.16%  "

Set-TextValue "D13" "11.11"
Set-TextValue "E13" "  +8.13%  "

Set-TextValue "D14" "4.687.71"
Set-TextValue "E14" "  +3.43%  "

Set-TextValue "D15" "4.060.58"
Set-TextValue "E15" "  +3.59%  "

Set-TextValue "D16" "21.17"
Set-TextValue "E16" "  +7.10%  "

Set-TextValue "D17" "14.14"
Set-TextValue "E17" "  +0.91%  "

Set-TextValue "E18" "  -1.29%  "

Set-TextValue "E19" "  -1.57%  "

Set-TextValue "D20" "72.342.67"
Set-TextValue "E20" "  +4.28%  "

Set-TextValue "D21" "445.37"
Set-TextValue "E21" "  +3.65%  "

Set-TextValue "D22" "104.89"
Set-TextValue "E22" "  +18.41%  "

Set-TextValue "D23" "3.62"
Set-TextValue "E23" "  +6.61%  "

Set-TextValue "D24" "14.76"
Set-TextValue "E24" "  +4.30%  "

Set-TextValue "E25" "  -1.41%  "

Set-TextValue "D26" "11.53"
Set-TextValue "E26" "  +0.32%  "

Set-TextValue "D27" "11.01"
Set-TextValue "E27" "  +3.46%  "

Set-TextValue "D28" "37.80"
Set-TextValue "E28" "  +3.60%  "

Set-TextValue "D29" "5.82"
Set-TextValue "E29" "  +2.62%  "

Set-TextValue "D30" "3.17"
Set-TextValue "E30" "  +11.80%  "

Set-TextValue "D31" "13.67"
Set-TextValue "E31" "  +3.50%  "

Set-TextValue "D32" "0.130"
Set-TextValue "E32" "  +2.97%  "

Set-TextValue "D33" "678.78"
Set-TextValue "E33" "  +0.10%  "

Set-TextValue "D34" "6.85"
Set-TextValue "E34" "  +15.09%  "

Set-TextValue "D35" "67.05"
Set-TextValue "E35" "  -2.86%  "

Set-TextValue "D36" "43.37"
Set-TextValue "E36" "  +8.24%  "

Set-TextValue "B37" "TheGraph"
Set-TextValue "C37" "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D37" "0.431"
Set-TextValue "E37" "  -1.84%  "

Set-TextValue "B38" "ThetaToken"
Set-TextValue "C38" "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue "D38" "3.62"
Set-TextValue "E38" "  +12.60%  "

Set-TextValue "D39" "0.0₃0859"
Set-TextValue "E39" "  -1.14%  "

Set-TextValue "E40" "  +1.45%  "

Set-TextValue "E41" "  +0.03%  "

Set-TextValue "D42" "0.0499"
Set-TextValue "E42" "  +3.79%  "

Set-TextValue "D43" "0.999"
Set-TextValue "E43" "  -0.17%  "

Set-TextValue "D44" "0.160"
Set-TextValue "E44" "  +13.45%  "

Set-TextValue "D45" "3.21"
Set-TextValue "E45" "  +1.45%  "

Set-TextValue "D46" "2.73"
Set-TextValue "E46" "  -2.43%  "

Set-TextValue "D47" "3.45"
Set-TextValue "E47" "  +3.70%  "

Set-TextValue "D48" "3.08"
Set-TextValue "E48" "  +2.60%  "

Set-TextValue "D49" "9.05"
Set-TextValue "E49" "  +7.22%  "

Set-TextValue "D50" "3.32"
Set-TextValue "E50" "  +1.80%  "

Set-TextValue "D51" "2.08"
Set-TextValue "E51" "  +0.64%  "

